$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 76, shifting existing rows 76-167 down to 77-168
$ws.Rows.Item(76).Insert()

# Populate the new row 76 with data
$ws.Range("A76").Value = 7
$ws.Range("B76").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C76").Value = "Ñuble"
$ws.Range("D76").Value = 44482
$ws.Range("E76").Value = 16
$ws.Range("F76").Value = 100112009
$ws.Range("G76").Value = "Acelga"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 120
$ws.Range("K76").Value = 300
$ws.Range("L76").Value = 350
$ws.Range("M76").Value = 325
$ws.Range("N76").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O76").Value = "Provincia de Diguillín"
$ws.Range("P76").Value = 325
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"
